$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("学習計画書")
Write-Host $ws.Name
